$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Death Line ("H") column for the three "Manage ..." sub-tasks under
# row 10 (rows 11-13) previously showed 2018/2019/2020 typos; correct them
# to the intended deadline "21h 15/1/2017" (matching row 10's deadline).
$ws.Range("H11").Value = "21h 15/1/2017"
$ws.Range("H12").Value = "21h 15/1/2017"
$ws.Range("H13").Value = "21h 15/1/2017"

# Reflect the author's last active selection when the file was saved.
$ws.Range("I12").Select()
